$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-11: apply the changed cell values (per the diff) ---
# Row 2
$ws.Range("E2").Value = 2
$ws.Range("G2").Value = 0.66889
$ws.Range("H2").Value = 1.33778
$ws.Range("I2").Value = 0.01180009979280774
$ws.Range("J2").Value = 0.01150792026947964
$ws.Range("K2").Value = 2
$ws.Range("M2").Value = 4.785814
$ws.Range("N2").Value = 9.571628
$ws.Range("O2").Value = 0.6750138823283496
$ws.Range("P2").Value = 0.5806605701408145
$ws.Range("Q2").Value = 3.20118312646
$ws.Range("R2").Value = 12.80473250584
$ws.Range("S2").Value = 0.007965231173005108
$ws.Range("T2").Value = 0.006682195544811084

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("G3").Value = 0.66889
$ws.Range("H3").Value = 1.33778
$ws.Range("I3").Value = 0.01180009979280774
$ws.Range("J3").Value = 0.01150792026947964
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.304135
$ws.Range("N3").Value = 6.912405000000001
$ws.Range("O3").Value = 0.3249861176716504
$ws.Range("P3").Value = 0.4193394298591856
$ws.Range("Q3").Value = 1.54121286015
$ws.Range("R3").Value = 9.247277160900001
$ws.Range("S3").Value = 0.003834868619802634
$ws.Range("T3").Value = 0.004825724724668559

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.349703666666667
$ws.Range("H4").Value = 7.049111
$ws.Range("I4").Value = 0.04145186465665943
$ws.Range("J4").Value = 0.0606382270318826
$ws.Range("K4").Value = 2
$ws.Range("M4").Value = 4.785814
$ws.Range("N4").Value = 9.571628
$ws.Range("O4").Value = 0.6750138823283496
$ws.Range("P4").Value = 0.5806605701408145
$ws.Range("Q4").Value = 11.24524470378467
$ws.Range("R4").Value = 67.471468222708
$ws.Range("S4").Value = 0.02798058409164098
$ws.Range("T4").Value = 0.0352102274806611

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.349703666666667
$ws.Range("H5").Value = 7.049111
$ws.Range("I5").Value = 0.04145186465665943
$ws.Range("J5").Value = 0.0606382270318826
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.304135
$ws.Range("N5").Value = 6.912405000000001
$ws.Range("O5").Value = 0.3249861176716504
$ws.Range("P5").Value = 0.4193394298591856
$ws.Range("Q5").Value = 5.414034457995
$ws.Range("R5").Value = 48.726310121955
$ws.Range("S5").Value = 0.01347128056501845
$ws.Range("T5").Value = 0.02542799955122151

# Row 6
$ws.Range("A6").Value = "M1"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.02107333333333333
$ws.Range("H6").Value = 0.06322
$ws.Range("I6").Value = 0.000371761330413723
$ws.Range("J6").Value = 0.0005438343520134124
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 4.785814
$ws.Range("N6").Value = 9.571628
$ws.Range("O6").Value = 0.6750138823283496
$ws.Range("P6").Value = 0.5806605701408145
$ws.Range("Q6").Value = 0.1008530536933333
$ws.Range("R6").Value = 0.60511832216
$ws.Range("S6").Value = 0.0002509440589421195
$ws.Range("T6").Value = 0.0003157831649022685

# Row 7
$ws.Range("A7").Value = "M1"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.02107333333333333
$ws.Range("H7").Value = 0.06322
$ws.Range("I7").Value = 0.000371761330413723
$ws.Range("J7").Value = 0.0005438343520134124
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 2.304135
$ws.Range("N7").Value = 6.912405000000001
$ws.Range("O7").Value = 0.3249861176716504
$ws.Range("P7").Value = 0.4193394298591856
$ws.Range("Q7").Value = 0.0485558049
$ws.Range("R7").Value = 0.4370022441
$ws.Range("S7").Value = 0.0001208172714716035
$ws.Range("T7").Value = 0.000228051187111144

# Row 8
$ws.Range("A8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.128441
$ws.Range("H8").Value = 0.385323
$ws.Range("I8").Value = 0.002265868255599604
$ws.Range("J8").Value = 0.00331464542899184
$ws.Range("K8").Value = 2
$ws.Range("M8").Value = 4.785814
$ws.Range("N8").Value = 9.571628
$ws.Range("O8").Value = 0.6750138823283496
$ws.Range("P8").Value = 0.5806605701408145
$ws.Range("Q8").Value = 0.614694735974
$ws.Range("R8").Value = 3.688168415844
$ws.Range("S8").Value = 0.001529492528056854
$ws.Range("T8").Value = 0.001924683904613047

# Row 9
$ws.Range("A9").Value = "M2"
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.128441
$ws.Range("H9").Value = 0.385323
$ws.Range("I9").Value = 0.002265868255599604
$ws.Range("J9").Value = 0.00331464542899184
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.304135
$ws.Range("N9").Value = 6.912405000000001
$ws.Range("O9").Value = 0.3249861176716504
$ws.Range("P9").Value = 0.4193394298591856
$ws.Range("Q9").Value = 0.295945403535
$ws.Range("R9").Value = 2.663508631815
$ws.Range("S9").Value = 0.0007363757275427502
$ws.Range("T9").Value = 0.001389961524378794

# Row 10
$ws.Range("A10").Value = "Neutro"
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.3791873333333333
$ws.Range("H10").Value = 1.137562
$ws.Range("I10").Value = 0.006689363532870856
$ws.Range("J10").Value = 0.009785594640067724
$ws.Range("K10").Value = 2
$ws.Range("M10").Value = 4.785814
$ws.Range("N10").Value = 9.571628
$ws.Range("O10").Value = 0.6750138823283496
$ws.Range("P10").Value = 0.5806605701408145
$ws.Range("Q10").Value = 1.814720048489333
$ws.Range("R10").Value = 10.888320290936
$ws.Range("S10").Value = 0.004515413248628841
$ws.Range("T10").Value = 0.005682108962868623

# Row 11
$ws.Range("A11").Value = "Neutro"
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.3791873333333333
$ws.Range("H11").Value = 1.137562
$ws.Range("I11").Value = 0.006689363532870856
$ws.Range("J11").Value = 0.009785594640067724
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 2.304135
$ws.Range("N11").Value = 6.912405000000001
$ws.Range("O11").Value = 0.3249861176716504
$ws.Range("P11").Value = 0.4193394298591856
$ws.Range("Q11").Value = 0.87369880629
$ws.Range("R11").Value = 7.863289256610001
$ws.Range("S11").Value = 0.002173950284242015
$ws.Range("T11").Value = 0.004103485677199102

# --- Insert 2 new rows (12-13) for the new "sCs" sending-cluster entries ---
$ws.Rows.Item(12).Resize(2).Insert() | Out-Null

# Row 12
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Ncam1"
$ws.Range("C12").Value = "Robo3"
$ws.Range("D12").Value = "ECs"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 53.1378185
$ws.Range("H12").Value = 106.275637
$ws.Range("I12").Value = 0.9374210424316487
$ws.Range("J12").Value = 0.9142097782775648
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 4.785814
$ws.Range("N12").Value = 9.571628
$ws.Range("O12").Value = 0.6750138823283496
$ws.Range("P12").Value = 0.5806605701408145
$ws.Range("Q12").Value = 254.307715706759
$ws.Range("R12").Value = 1017.230862827036
$ws.Range("S12").Value = 0.6327722172280758
$ws.Range("T12").Value = 0.5308455710829584

# Row 13
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Ncam1"
$ws.Range("C13").Value = "Robo3"
$ws.Range("D13").Value = "Neutro"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 53.1378185
$ws.Range("H13").Value = 106.275637
$ws.Range("I13").Value = 0.9374210424316487
$ws.Range("J13").Value = 0.9142097782775648
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 2.304135
$ws.Range("N13").Value = 6.912405000000001
$ws.Range("O13").Value = 0.3249861176716504
$ws.Range("P13").Value = 0.4193394298591856
$ws.Range("Q13").Value = 122.4367074294975
$ws.Range("R13").Value = 734.6202445769851
$ws.Range("S13").Value = 0.3046488252035729
$ws.Range("T13").Value = 0.3833642071946065

Write-Output "Applied 148 cell updates + inserted rows 12-13"
